# Generate Report for Handoff
#
# - Flip every "In Translation" status cell to "Ready for handoff"
#   (Overview!E2:F2, zh-cn!C2, de-de!C2).
# - Bump the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#   timestamps to reflect the new handoff-ready report.
# - Widen the Status columns so the longer "Ready for handoff" text fits
#   (mirrors Excel's column-width auto-adjustment).

$wb = $excel.ActiveWorkbook

# --- Update status text on every sheet ------------------------------------
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("In Translation", "Ready for handoff")
}

# --- Update timestamps -----------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-28 02:39:08"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-28 02:39:01"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-28 02:39:08"

# --- Widen the Status columns to fit the new, longer text ------------------
$newStatusColumnWidth = 16.3333333333333

$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColumnWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColumnWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColumnWidth
